$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cellRef, $val)
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.ClearFormats()
}

Set-TextValue 'D2' '62.699.64'
Set-TextValue 'E2' '  +6.20%  '
Set-TextValue 'D3' '3.108.01'
Set-TextValue 'E3' '  +3.66%  '
Set-TextValue 'E4' '  +0.20%  '
Set-TextValue 'D5' '584.96'
Set-TextValue 'D6' '143.56'
Set-TextValue 'E6' '  +4.98%  '
Set-TextValue 'E7' '  +0.01%  '
Set-TextValue 'D8' '3.097.05'
Set-TextValue 'E8' '  +3.50%  '
Set-TextValue 'E9' '  +1.87%  '
Set-TextValue 'E10' '  +8.84%  '
Set-TextValue 'D11' '5.74'
Set-TextValue 'E11' '  +9.66%  '
Set-TextValue 'E12' '  +2.75%  '
Set-TextValue 'D13' '0.0000244'
Set-TextValue 'E13' '  +5.81%  '
Set-TextValue 'D14' '35.58'
Set-TextValue 'E14' '  +5.92%  '
Set-TextValue 'E15' '  +0.84%  '
Set-TextValue 'D16' '3.624.46'
Set-TextValue 'E16' '  +3.89%  '
Set-TextValue 'D17' '7.28'
Set-TextValue 'E17' '  -0.08%  '
Set-TextValue 'D18' '3.107.49'
Set-TextValue 'E18' '  +3.72%  '
Set-TextValue 'D19' '62.680.29'
Set-TextValue 'E19' '  +6.15%  '
Set-TextValue 'D20' '453.40'
Set-TextValue 'E20' '  +5.76%  '
Set-TextValue 'D21' '14.07'
Set-TextValue 'E21' '  +2.72%  '
Set-TextValue 'D22' '0.734'
Set-TextValue 'E22' '  +1.88%  '
Set-TextValue 'E23' '  +6.05%  '
Set-TextValue 'D24' '13.69'
Set-TextValue 'E24' '  +2.66%  '
Set-TextValue 'D25' '82.09'
Set-TextValue 'E25' '  +1.78%  '
Set-TextValue 'E27' '  +4.21%  '
Set-TextValue 'E28' '  +5.90%  '
Set-TextValue 'E29' '  +0.21%  '
Set-TextValue 'E30' '  +5.50%  '
Set-TextValue 'E31' '  +13.96%  '
Set-TextValue 'E32' '  +12.26%  '
Set-TextValue 'D33' '27.13'
Set-TextValue 'E34' '  +5.06%  '
Set-TextValue 'D35' '0.0₃0801'
Set-TextValue 'E35' '  +5.11%  '
Set-TextValue 'D36' '6.08'
Set-TextValue 'E36' '  +2.29%  '
Set-TextValue 'D37' '2.21'
Set-TextValue 'E37' '  +5.56%  '
Set-TextValue 'D38' '50.65'
Set-TextValue 'E38' '  +3.88%  '
Set-TextValue 'E39' '  +10.98%  '
Set-TextValue 'E40' '  +1.56%  '
Set-TextValue 'D41' '425.16'
Set-TextValue 'E41' '  +6.49%  '
Set-TextValue 'D42' '2.932.86'
Set-TextValue 'E42' '  +6.40%  '
Set-TextValue 'E43' '  +6.06%  '
Set-TextValue 'E44' '  +12.48%  '
Set-TextValue 'E45' '  +3.69%  '
Set-TextValue 'E46' '  +8.66%  '
Set-TextValue 'D47' '125.30'
Set-TextValue 'E47' '  +1.50%  '
Set-TextValue 'D49' '34.89'
Set-TextValue 'E49' '  -1.96%  '
Set-TextValue 'E50' '  +1.28%  '
Set-TextValue 'D51' '24.76'
Set-TextValue 'E51' '  +5.82%  '
